$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 22.77000000000012
$ws.Range("G2").Value = [double]"1.110223024625157e-16"
$ws.Range("H2").Value = [double]"2.254259948477475e-16"
$ws.Range("K2").Value = 48.51248143355631
$ws.Range("L2").Value = "[40.53957669627214, 56.48538617084048]"
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 1.540921321580579
$ws.Range("P2").Value = "[1.3522370781217319, 1.7296055650394262]"
$ws.Range("S2").Value = 65.09931036773861
$ws.Range("T2").Value = "[59.71206483042523, 70.486555905052]"
$ws.Range("W2").Value = 17.18576576576586
$ws.Range("X2").Value = 16.50198198198207
$ws.Range("Y2").Value = 17.86954954954965

# Row 3 updates
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 25.56000000000056
$ws.Range("G3").Value = [double]"1.110223024625157e-16"
$ws.Range("H3").Value = [double]"2.254259948477475e-16"
$ws.Range("I3").Value = ""
$ws.Range("K3").Value = 44.5362795641766
$ws.Range("L3").Value = "[33.3285599290011, 55.7439991993521]"
$ws.Range("M3").Value = [double]"1.021405182655144e-13"
$ws.Range("N3").Value = [double]"1.021405182655144e-13"
$ws.Range("O3").Value = -0.6415264277600778
$ws.Range("P3").Value = "[-0.8805264694746171, -0.4025263860455386]"
$ws.Range("Q3").Value = [double]"2.523152249356286e-07"
$ws.Range("R3").Value = [double]"2.523152249356286e-07"
$ws.Range("S3").Value = 63.27589039532914
$ws.Range("T3").Value = "[57.382677073339416, 69.16910371731888]"
$ws.Range("W3").Value = 2.609729729729786
$ws.Range("X3").Value = 1.63747747747751
$ws.Range("Y3").Value = 3.581981981982061
